$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/centered header style used by A1:AC1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the team's win/loss/tie record for every data row (2..46)
for ($row = 2; $row -le 46; $row++) {
    $ws.Cells.Item($row, 30).Value = 78  # AD: Wins
    $ws.Cells.Item($row, 31).Value = 83  # AE: Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF: Ties
}
